$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "257.89"
Set-TextValue $ws.Range("E2") "1.45%"
Set-TextValue $ws.Range("D3") "26.85"
Set-TextValue $ws.Range("E3") "-4.06%"
Set-TextValue $ws.Range("D4") "4.736"
Set-TextValue $ws.Range("E4") "-10.96%"
Set-TextValue $ws.Range("D5") "0.05977"
Set-TextValue $ws.Range("E5") "2.32%"
Set-TextValue $ws.Range("D6") "6.686"
Set-TextValue $ws.Range("E6") "-0.22%"
Set-TextValue $ws.Range("D7") "0.8740"
Set-TextValue $ws.Range("E7") "0.85%"
Set-TextValue $ws.Range("D8") "0.9550"
Set-TextValue $ws.Range("E8") "4.42%"
Set-TextValue $ws.Range("D9") "0.1414"
Set-TextValue $ws.Range("E9") "-1.01%"
Set-TextValue $ws.Range("D10") "0.07243"
Set-TextValue $ws.Range("E10") "0.88%"
Set-TextValue $ws.Range("D11") "0.03138"
Set-TextValue $ws.Range("E11") "-2.29%"
Set-TextValue $ws.Range("D12") "0.09239"
Set-TextValue $ws.Range("E12") "0.14%"
Set-TextValue $ws.Range("D13") "0.001555"
Set-TextValue $ws.Range("E13") "0.15%"
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D14") "0.0006121"
Set-TextValue $ws.Range("E14") "0.55%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D15") "0.006021"
Set-TextValue $ws.Range("E15") "0.51%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D16") "3.487"
Set-TextValue $ws.Range("E16") "-0.30%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D17") "3.205"
Set-TextValue $ws.Range("E17") "-0.56%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D18") "2.219"
Set-TextValue $ws.Range("E18") "-1.40%"
Set-TextValue $ws.Range("D19") "0.3114"
Set-TextValue $ws.Range("E19") "-1.69%"
Set-TextValue $ws.Range("D20") "0.03619"
Set-TextValue $ws.Range("E20") "4.84%"
Set-TextValue $ws.Range("E21") "-1.47%"
Set-TextValue $ws.Range("D22") "3.537"
Set-TextValue $ws.Range("E22") "0.23%"
Set-TextValue $ws.Range("D23") "0.04230"
Set-TextValue $ws.Range("E23") "2.11%"
Set-TextValue $ws.Range("E24") "2.63%"
Set-TextValue $ws.Range("D25") "0.001214"
Set-TextValue $ws.Range("E25") "-1.03%"
Set-TextValue $ws.Range("E26") "-11.74%"
Set-TextValue $ws.Range("D27") "0.0001201"
Set-TextValue $ws.Range("E27") "0.11%"
Set-TextValue $ws.Range("E28") "-22.92%"
Set-TextValue $ws.Range("D40") "0.03833"
Set-TextValue $ws.Range("E40") "-0.38%"
Set-TextValue $ws.Range("D41") "0.006184"
Set-TextValue $ws.Range("E41") "9.12%"
Set-TextValue $ws.Range("E42") "0.37%"
Set-TextValue $ws.Range("D43") "0.002301"
Set-TextValue $ws.Range("E43") "4.66%"
Set-TextValue $ws.Range("D44") "0.01053"
Set-TextValue $ws.Range("E44") "-4.38%"
Set-TextValue $ws.Range("D45") "0.00005496"
Set-TextValue $ws.Range("E45") "4.16%"
Set-TextValue $ws.Range("E46") "0.06%"
Set-TextValue $ws.Range("D47") "0.08551"
Set-TextValue $ws.Range("E47") "-21.40%"
Set-TextValue $ws.Range("D48") "0.002132"
Set-TextValue $ws.Range("E48") "-1.75%"
Set-TextValue $ws.Range("E49") "0.06%"
Set-TextValue $ws.Range("E50") "0.06%"
